$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.614.04'
$ws.Range("E2").Value = '  -3.72%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.864.56'
$ws.Range("E3").Value = '  -4.22%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '509.53'
$ws.Range("E5").Value = '  -6.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.38'
$ws.Range("E6").Value = '  -8.75%  '
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("E8").Value = '  -6.69%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.859.29'
$ws.Range("E9").Value = '  -4.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.06'
$ws.Range("E10").Value = '  -1.21%  '
$ws.Range("E11").Value = '  -8.80%  '
$ws.Range("E12").Value = '  -5.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.363.06'
$ws.Range("E13").Value = '  -4.45%  '
$ws.Range("E14").Value = '  +1.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '59.944.71'
$ws.Range("E15").Value = '  -3.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.90'
$ws.Range("E16").Value = '  -8.50%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.855.53'
$ws.Range("E17").Value = '  -4.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000137'
$ws.Range("E18").Value = '  -7.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.78'
$ws.Range("E19").Value = '  -7.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.20'
$ws.Range("E20").Value = '  -6.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '350.67'
$ws.Range("E21").Value = '  -7.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.34'
$ws.Range("E22").Value = '  -5.67%  '
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("E24").Value = '  -0.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.26'
$ws.Range("E25").Value = '  -4.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.436'
$ws.Range("E26").Value = '  -7.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.173'
$ws.Range("E27").Value = '  -8.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.56%  '
$ws.Range("E29").Value = '  -8.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0833'
$ws.Range("E30").Value = '  -10.71%  '
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.12'
$ws.Range("E33").Value = '  -6.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '151.78'
$ws.Range("E34").Value = '  -5.70%  '
$ws.Range("E35").Value = '  -8.14%  '
$ws.Range("E36").Value = '  -8.46%  '
$ws.Range("E37").Value = '  -10.29%  '
$ws.Range("E38").Value = '  -8.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.96'
$ws.Range("E39").Value = '  -1.56%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.238.84'
$ws.Range("E40").Value = '  -7.39%  '
$ws.Range("B41").Value = 'Mantle'
$ws.Range("C41").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.637'
$ws.Range("E41").Value = '  -5.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.41'
$ws.Range("E42").Value = '  -9.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.56'
$ws.Range("E43").Value = '  -8.82%  '
$ws.Range("E44").Value = '  -4.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.81'
$ws.Range("E46").Value = '  -10.16%  '
$ws.Range("E47").Value = '  -0.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0229'
$ws.Range("E48").Value = '  -6.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.67'
$ws.Range("E49").Value = '  -12.70%  '
$ws.Range("E50").Value = '  -6.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.80'
$ws.Range("E51").Value = '  -9.00%  '
